# Add team record (Wins/Losses/Ties) columns to the sheet,
# matching the existing header style used by the other header cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (style) from an existing header cell (A1)
# onto the new header cells before writing their text.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AC1:AE1").PasteSpecial(-4122) | Out-Null # xlPasteFormats

# Header row (row 1) - new columns AC, AD, AE
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Data rows (2-42) - team win/loss/tie record for every player row
$lastRow = 42
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 29).Value = 97  # AC - Wins
    $ws.Cells.Item($r, 30).Value = 65  # AD - Losses
    $ws.Cells.Item($r, 31).Value = 0   # AE - Ties
}
